$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Werkstoffeigenschaften")
$rng = $ws2.Range("B9:D9")
$rng.Borders.LineStyle = 1
$rng.Borders.Weight = 2
$rng.Borders.Color = RGB(166,166,166)
Write-Output "done"
